$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 65 into row 66 first. This preserves the original text-typed
# "3" that was in B65 so the new B66 keeps the same (string) representation,
# matching how the source data actually evolved.
$ws.Range("A65:H65").Copy($ws.Range("A66:H66"))

# B65 itself is retyped as a genuine number (3) rather than the text "3" it
# held before.
$ws.Range("B65").Value = 3

# Update the fields on the newly appended row 66 with the new annotation.
$ws.Range("D66").Value = "DFT"
$ws.Range("E66").Value = "WRI"
$ws.Range("F66").Value = "c1f109d6-e04a-469c-a254-426c0826b7a8"
$ws.Range("G66").Value = "BJcAWaeCW_annotated.xlsx"
$ws.Range("H66").Value = "The paper organization needs work; there are also some missing pieces to put the NN training together."
